$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; this pushes the existing rows 15-41
# down to 16-42 (matching the diff, which shows every record from the old
# row 15 onward shifted down by one row, with the old row 41 duplicated
# into the new row 42 untouched, and a brand-new record occupying row 15).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(15, 3).Value = "Maule"
$ws.Cells.Item(15, 4).Value = 44757
$ws.Cells.Item(15, 5).Value = 7
$ws.Cells.Item(15, 6).Value = 100112043
$ws.Cells.Item(15, 7).Value = "Pepino dulce"
$ws.Cells.Item(15, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 15000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 15000
$ws.Cells.Item(15, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 833
$ws.Cells.Item(15, 17).Value = 18
$ws.Cells.Item(15, 18).Value = "Hortaliza"
